$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update corrected values on existing rows 211 and 213 ---
$ws.Range("C211:F211").Value = 9409316708800
$ws.Range("C213:F213").Value = 11478475285200

# --- Append new rows 214-216 ---
$ws.Range("A214").Value = 45139.41666666666
$ws.Range("B214").Value = "ECONOMICS:TRM2"
$ws.Range("C214:F214").Value = 12025285811100
$ws.Range("G214").Value = 0

$ws.Range("A215").Value = 45170.41666666666
$ws.Range("B215").Value = "ECONOMICS:TRM2"
$ws.Range("C215:F215").Value = 12349311426800
$ws.Range("G215").Value = 0

$ws.Range("A216").Value = 45200.45833333334
$ws.Range("B216").Value = "ECONOMICS:TRM2"
$ws.Range("C216:F216").Value = 12763732048500
$ws.Range("G216").Value = 0

# Copy style (date format / bold / border) from the row above into the new date cells
$ws.Range("A213").Copy()
$ws.Range("A214:A216").PasteSpecial(-4122) | Out-Null
